$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "A4"   = -20.61979999999998
    "A7"   = -22.02250000000002
    "B7"   = 4.626900000000001
    "B15"  = 4.885099999999995
    "A16"  = -21.62999999999998
    "C16"  = -11.83739999999999
    "C19"  = -12.00070000000001
    "B21"  = 9.989400000000003
    "B22"  = 10.17610000000001
    "B23"  = 8.985300000000009
    "A28"  = -22.03679999999999
    "A29"  = -21.06629999999997
    "A32"  = -21.21449999999999
    "B34"  = 9.511600000000005
    "E34"  = 17.40760000000001
    "C36"  = -12.6761
    "A40"  = -19.8768
    "B43"  = 6.127000000000002
    "E43"  = 17.37390000000002
    "B45"  = 4.810500000000003
    "C46"  = -14.69
    "E48"  = 17.467
    "B50"  = 4.706699999999999
    "C50"  = -13.49539999999999
    "B51"  = 5.7702
    "A52"  = -22.17209999999999
    "A57"  = -22.49110000000001
    "A66"  = -21.4817
    "B66"  = 4.980399999999998
    "B67"  = 5.0076
    "E70"  = 17.61140000000001
    "E73"  = 16.93930000000001
    "B79"  = 9.583500000000008
    "B84"  = 5.363
    "E87"  = 16.37999999999999
    "B92"  = 4.721699999999999
    "E92"  = 18.55100000000001
    "C95"  = -12.7495
    "B97"  = 6.553599999999997
    "C97"  = -12.3224
    "A100" = -22.06199999999999
    "E101" = 16.99910000000001
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
